$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LUY_KE_THANG_LONG_XUYEN")

# Update last_edited_time (column D) for rows 3, 4, 5, 7, 13
$ws.Range("D3").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("D4").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("D5").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("D7").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("D13").Value = "2024-07-18T15:58:00.000Z"

# Update numeric values on row 5 (Chi tieu / Luy ke formulas)
$ws.Range("W5").Value = 16617000
$ws.Range("AA5").Value = 21333000
